# Slide 13 ("If you rolled 3000 times / and did not get 425-525 / ...")
# The number callout "425" is corrected to "475" while the "-525" portion
# is kept (as its own run, same formatting) immediately after it.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(13)

foreach ($shape in $s.Shapes) {
    if ($shape.Name -eq "TextBox 3") {
        $tr2 = $shape.TextFrame2.TextRange
        $fullText = $tr2.Text
        $idx = $fullText.IndexOf("425") + 1   # 1-based COM character index
        $chars = $tr2.Characters($idx, 3)
        $chars.Text = "475"
    }
}
